$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column I (RF) for rows 23 through 53 from 5.6966 to 2.947435897435897
$ws.Range("I23:I53").Value = 2.947435897435897
